# issue #5: stock data output to json file
# Add a "property_category" column (value "stock") to the 股票 (stock) sheet,
# inserted right after the "total" column and before the "date" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column at H (shifts old H:J -> I:K, i.e. date / legislator_name /
# legislator_id move one column to the right) so the existing formatting of the
# neighbouring cells carries over to the shifted cells.
$ws.Columns.Item(8).Insert()

# Header for the freshly inserted column, formatted like the other header cells.
$ws.Cells.Item(1, 8).Value = "property_category"
$ws.Cells.Item(1, 8).Font.Bold = $true
$ws.Cells.Item(1, 8).Borders.LineStyle = 1
$ws.Cells.Item(1, 8).HorizontalAlignment = -4108
$ws.Cells.Item(1, 8).VerticalAlignment = -4160

# Every stock row is a "stock" category property.
$lastRow = $ws.UsedRange.Row + $ws.UsedRange.Rows.Count - 1
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 8).Value = "stock"
}
